# Insert a new data row at row 5 (pushing the existing rows 5-84 down to 6-85)
# and populate it with a new "Fuyu / Primera" observation dated 45092 with a
# volume of 100, mirroring the other fields of the record that used to sit
# at row 5 (which now lives at row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 45092
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107001
$ws.Range("J5").Value = "Caqui"
$ws.Range("K5").Value = "Fuyu"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("Q5").Value = "$/bandeja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1600
$ws.Range("T5").Value = 15
